# Apply the "Updated symbol list" data refresh to the cryptos worksheet.
# All target cells hold text-formatted values (prices, percentages, coin names,
# and links) stored as inline strings, so we force NumberFormat "@" (Text) before
# assigning values in the numeric-looking Price/Volume columns (D, E) to prevent
# Excel from auto-converting them to numbers and losing exact text formatting
# (trailing zeros, percent signs, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '305.63'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.88%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.67'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.60%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.061'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.41%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08113'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.87%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.969'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.08%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.170'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.00%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '7.762'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.50%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9308'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.61%'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '10.67%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1908'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.80%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09273'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.76%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03551'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2.75%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09860'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.06%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001412'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.45%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005756'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.12%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.560'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.49%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.045'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.82%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3443'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.23%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1344'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '4.79%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.893'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-3.16%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2595'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '8.26%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04402'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.89%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001219'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.76%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004785'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.55%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '32.19%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003125'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '4.26%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01969'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2.42%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04966'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '4.71%'
$ws.Range("B41").Value = 'Dexo'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QkL_pl546+dexo-dexo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01083'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '12.31%'
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007610'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.01%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1380'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002098'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.46%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01080'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.31%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006384'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2.02%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.01%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001189'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-20.16%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002098'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.01%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001998'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.01%'
